$d = $word.ActiveDocument

# Step 1: delete paragraphs 69 through 81 inclusive (the empty centered
# paragraph after the "Diagrama de Secuencia" picture, through the last
# "Diagrama de Estados" comment paragraphs including the trailing blank
# PSI-Comentario paragraph right before the final PSI-Ttulo1 paragraph).
$startPara = $d.Paragraphs.Item(69)
$endPara = $d.Paragraphs.Item(81)
$delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$delRange.Delete()

# Step 2: remove the picture ("Sin título-5.jpg") that lives inside the
# final paragraph, leaving the paragraph itself (and its PSI-Ttulo1 style)
# intact but empty.
$lastShape = $d.InlineShapes.Item($d.InlineShapes.Count)
$lastShape.Delete()
